# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (Castle Brite / Especial, Primera, Segunda)
# ahead of the existing row 13, pushing the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 3 new records right above what is currently row 13.
$ws.Range("A13:T15").Insert()

# --- New row 13: Castle Brite / Especial ---
$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44540
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100103
$ws.Range("H13").Value = "Frutos de hueso (carozo)"
$ws.Range("I13").Value = 100103003
$ws.Range("J13").Value = "Damasco"
$ws.Range("K13").Value = "Castle Brite"
$ws.Range("L13").Value = "Especial"
$ws.Range("M13").Value = 140
$ws.Range("N13").Value = 24500
$ws.Range("O13").Value = 25000
$ws.Range("P13").Value = 24750
$ws.Range("Q13").Value = "$/caja 15 kilos"
$ws.Range("R13").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S13").Value = 1650
$ws.Range("T13").Value = 15

# --- New row 14: Castle Brite / Primera ---
$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "Terminal La Palmera de La Serena"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44540
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100103
$ws.Range("H14").Value = "Frutos de hueso (carozo)"
$ws.Range("I14").Value = 100103003
$ws.Range("J14").Value = "Damasco"
$ws.Range("K14").Value = "Castle Brite"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 22500
$ws.Range("O14").Value = 23000
$ws.Range("P14").Value = 22750
$ws.Range("Q14").Value = "$/caja 15 kilos"
$ws.Range("R14").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S14").Value = 1517
$ws.Range("T14").Value = 15

# --- New row 15: Castle Brite / Segunda ---
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = "Terminal La Palmera de La Serena"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44540
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100103
$ws.Range("H15").Value = "Frutos de hueso (carozo)"
$ws.Range("I15").Value = 100103003
$ws.Range("J15").Value = "Damasco"
$ws.Range("K15").Value = "Castle Brite"
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 18000
$ws.Range("O15").Value = 18500
$ws.Range("P15").Value = 18250
$ws.Range("Q15").Value = "$/caja 15 kilos"
$ws.Range("R15").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S15").Value = 1217
$ws.Range("T15").Value = 15
